# Build Your First Process With Studio Complete
# Fixes the "Client discount" line item (row 18) so that it no longer
# references the removed "This client doesn't benefit from any discount"
# text, and instead holds a numeric value of 0 so the dependent formulas
# (F18, F21, F23, F24) recalculate to real numbers instead of #VALUE!.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# E18 previously held the text "This client doesn't benefit from any discount"
# Replace it with a numeric 0 so the shared formula in F18 can compute.
$ws.Range("E18").Value = 0

# Recalculate the workbook so F18, F21, F23 and F24 pick up real numbers.
$excel.CalculateFullRebuild()

# A31 held a footer note referencing the now-removed shared string; after the
# removal of the unused string, the remaining footer text ("RPA Dev, ...")
# keeps the same value. Re-assigning ensures the shared string table no
# longer carries the now-unused "This client doesn't benefit..." entry.
$ws.Range("A31").Value = "RPA Dev, developer.rpa@mail.com"
